$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 odds updates
$ws.Range("G2").Value = 2.9
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 2.88
$ws.Range("L2").Value = 4
$ws.Range("N2").Value = 5
$ws.Range("AH2").Value = 12
$ws.Range("AJ2").Value = 34
$ws.Range("AW2").Value = 4.5

# Row 4 odds updates
$ws.Range("G4").Value = 2.25
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 3.1
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("AQ4").Value = 51
$ws.Range("AZ4").Value = 67
$ws.Range("BA4").Value = 101
$ws.Range("BB4").Value = 301

# Remove the last data row (row 6) - Paraguay match no longer present
$ws.Rows.Item(6).Delete()
